$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("1.0", "Automation Test Engineer (UFT)", "https://www.naukri.com/job-listings-automation-test-engineer-uft-tata-consultancy-services-chennai-bengaluru-5-to-10-years-080126026105", "3+ weeks ago", "Apply"),
    @("2.0", "Automation Testing Engineer- Playwright", "https://www.naukri.com/job-listings-automation-testing-engineer-playwright-sonata-software-hyderabad-chennai-bengaluru-7-to-12-years-231225035819", "3+ weeks ago", "Apply"),
    @("3.0", "Automation Test Engineer(C#+API)", "https://www.naukri.com/job-listings-automation-test-engineer-c-api-ltimindtree-hyderabad-chennai-bengaluru-5-to-10-years-020126008100", "3+ weeks ago", "Apply"),
    @("4.0", "Test Automation Engineer", "https://www.naukri.com/job-listings-test-automation-engineer-accenture-solutions-pvt-ltd-chennai-3-to-8-years-130126921557", "3+ weeks ago", "ApplyOnCompanySite"),
    @("5.0", "Storage Test Automation Engineer", "https://www.naukri.com/job-listings-storage-test-automation-engineer-aziro-pune-chennai-bengaluru-4-to-9-years-110126008534", "3+ weeks ago", "Apply"),
    @("6.0", "Cypress Automation Test Engineer", "https://www.naukri.com/job-listings-cypress-automation-test-engineer-cognizant-chennai-4-to-10-years-020126013272", "3+ weeks ago", "Apply"),
    @("7.0", "Test Automation Engineer | Contract", "https://www.naukri.com/job-listings-test-automation-engineer-contract-srs-infoway-chennai-8-to-10-years-090126939555", "3+ weeks ago", "Apply"),
    @("8.0", "IVI Test Automation Engineer (Automotive Infotainment)", "https://www.naukri.com/job-listings-ivi-test-automation-engineer-automotive-infotainment-tekpillar-pune-chennai-bengaluru-5-to-10-years-100126000791", "3+ weeks ago", "Apply"),
    @("9.0", "IVR Automation Test Engineer(Cyara & Kore.ai)", "https://www.naukri.com/job-listings-ivr-automation-test-engineer-cyara-kore-ai-summitworks-technologies-pune-chennai-5-to-10-years-030126013225", "3+ weeks ago", "Apply"),
    @("10.0", "Automation Test Engineer", "https://www.naukri.com/job-listings-automation-test-engineer-indium-software-chennai-bengaluru-5-to-10-years-090126019771", "3+ weeks ago", "Apply"),
    @("11.0", "Automation Test Engineer (UFT)", "https://www.naukri.com/job-listings-automation-test-engineer-uft-tata-consultancy-services-chennai-bengaluru-5-to-10-years-080126026105", "Today", "Apply"),
    @("12.0", "Automation Testing Engineer- Playwright", "https://www.naukri.com/job-listings-automation-testing-engineer-playwright-sonata-software-hyderabad-chennai-bengaluru-7-to-12-years-231225035819", "Today", "Apply"),
    @("13.0", "Automation Test Engineer(C#+API)", "https://www.naukri.com/job-listings-automation-test-engineer-c-api-ltimindtree-hyderabad-chennai-bengaluru-5-to-10-years-020126008100", "Today", "Apply"),
    @("14.0", "Test Automation Engineer", "https://www.naukri.com/job-listings-test-automation-engineer-accenture-solutions-pvt-ltd-chennai-3-to-8-years-130126921557", "Today", "ApplyOnCompanySite"),
    @("15.0", "Storage Test Automation Engineer", "https://www.naukri.com/job-listings-storage-test-automation-engineer-aziro-pune-chennai-bengaluru-4-to-9-years-110126008534", "Today", "Apply"),
    @("16.0", "Cypress Automation Test Engineer", "https://www.naukri.com/job-listings-cypress-automation-test-engineer-cognizant-chennai-4-to-10-years-020126013272", "Today", "Apply"),
    @("17.0", "Test Automation Engineer | Contract", "https://www.naukri.com/job-listings-test-automation-engineer-contract-srs-infoway-chennai-8-to-10-years-090126939555", "Today", "Apply"),
    @("18.0", "IVI Test Automation Engineer (Automotive Infotainment)", "https://www.naukri.com/job-listings-ivi-test-automation-engineer-automotive-infotainment-tekpillar-pune-chennai-bengaluru-5-to-10-years-100126000791", "Today", "Apply"),
    @("19.0", "IVR Automation Test Engineer(Cyara & Kore.ai)", "https://www.naukri.com/job-listings-ivr-automation-test-engineer-cyara-kore-ai-summitworks-technologies-pune-chennai-5-to-10-years-030126013225", "Today", "Apply"),
    @("20.0", "Automation Test Engineer", "https://www.naukri.com/job-listings-automation-test-engineer-indium-software-chennai-bengaluru-5-to-10-years-090126019771", "Today", "Apply"),
    @("21.0", "Automation Test Engineer (UFT)", "https://www.naukri.com/job-listings-automation-test-engineer-uft-tata-consultancy-services-chennai-bengaluru-5-to-10-years-080126026105", "4 days ago", "Apply"),
    @("22.0", "Automation Testing Engineer- Playwright", "https://www.naukri.com/job-listings-automation-testing-engineer-playwright-sonata-software-hyderabad-chennai-bengaluru-7-to-12-years-231225035819", "4 days ago", "Apply"),
    @("23.0", "Automation Test Engineer(C#+API)", "https://www.naukri.com/job-listings-automation-test-engineer-c-api-ltimindtree-hyderabad-chennai-bengaluru-5-to-10-years-020126008100", "4 days ago", "Apply"),
    @("24.0", "Test Automation Engineer", "https://www.naukri.com/job-listings-test-automation-engineer-accenture-solutions-pvt-ltd-chennai-3-to-8-years-130126921557", "4 days ago", "ApplyOnCompanySite"),
    @("25.0", "Storage Test Automation Engineer", "https://www.naukri.com/job-listings-storage-test-automation-engineer-aziro-pune-chennai-bengaluru-4-to-9-years-110126008534", "4 days ago", "Apply"),
    @("26.0", "Cypress Automation Test Engineer", "https://www.naukri.com/job-listings-cypress-automation-test-engineer-cognizant-chennai-4-to-10-years-020126013272", "4 days ago", "Apply"),
    @("27.0", "Test Automation Engineer | Contract", "https://www.naukri.com/job-listings-test-automation-engineer-contract-srs-infoway-chennai-8-to-10-years-090126939555", "4 days ago", "Apply"),
    @("28.0", "IVI Test Automation Engineer (Automotive Infotainment)", "https://www.naukri.com/job-listings-ivi-test-automation-engineer-automotive-infotainment-tekpillar-pune-chennai-bengaluru-5-to-10-years-100126000791", "4 days ago", "Apply"),
    @("29.0", "IVR Automation Test Engineer(Cyara & Kore.ai)", "https://www.naukri.com/job-listings-ivr-automation-test-engineer-cyara-kore-ai-summitworks-technologies-pune-chennai-5-to-10-years-030126013225", "4 days ago", "Apply"),
    @("30.0", "Automation Test Engineer", "https://www.naukri.com/job-listings-automation-test-engineer-indium-software-chennai-bengaluru-5-to-10-years-090126019771", "4 days ago", "Apply"),
    @("31.0", "Automation Test Engineer (UFT)", "https://www.naukri.com/job-listings-automation-test-engineer-uft-tata-consultancy-services-chennai-bengaluru-5-to-10-years-080126026105", "1 week ago", "Apply"),
    @("32.0", "Automation Testing Engineer- Playwright", "https://www.naukri.com/job-listings-automation-testing-engineer-playwright-sonata-software-hyderabad-chennai-bengaluru-7-to-12-years-231225035819", "1 week ago", "Apply"),
    @("33.0", "Automation Test Engineer(C#+API)", "https://www.naukri.com/job-listings-automation-test-engineer-c-api-ltimindtree-hyderabad-chennai-bengaluru-5-to-10-years-020126008100", "1 week ago", "Apply"),
    @("34.0", "Test Automation Engineer", "https://www.naukri.com/job-listings-test-automation-engineer-accenture-solutions-pvt-ltd-chennai-3-to-8-years-130126921557", "1 week ago", "ApplyOnCompanySite"),
    @("35.0", "Storage Test Automation Engineer", "https://www.naukri.com/job-listings-storage-test-automation-engineer-aziro-pune-chennai-bengaluru-4-to-9-years-110126008534", "1 week ago", "Apply"),
    @("36.0", "Cypress Automation Test Engineer", "https://www.naukri.com/job-listings-cypress-automation-test-engineer-cognizant-chennai-4-to-10-years-020126013272", "1 week ago", "Apply"),
    @("37.0", "Test Automation Engineer | Contract", "https://www.naukri.com/job-listings-test-automation-engineer-contract-srs-infoway-chennai-8-to-10-years-090126939555", "1 week ago", "Apply"),
    @("38.0", "IVI Test Automation Engineer (Automotive Infotainment)", "https://www.naukri.com/job-listings-ivi-test-automation-engineer-automotive-infotainment-tekpillar-pune-chennai-bengaluru-5-to-10-years-100126000791", "1 week ago", "Apply"),
    @("39.0", "IVR Automation Test Engineer(Cyara & Kore.ai)", "https://www.naukri.com/job-listings-ivr-automation-test-engineer-cyara-kore-ai-summitworks-technologies-pune-chennai-5-to-10-years-030126013225", "1 week ago", "Apply"),
    @("40.0", "Automation Test Engineer", "https://www.naukri.com/job-listings-automation-test-engineer-indium-software-chennai-bengaluru-5-to-10-years-090126019771", "1 week ago", "Apply"),
    @("41.0", "Automation Test Engineer (UFT)", "https://www.naukri.com/job-listings-automation-test-engineer-uft-tata-consultancy-services-chennai-bengaluru-5-to-10-years-080126026105", "5 days ago", "Apply"),
    @("42.0", "Automation Testing Engineer- Playwright", "https://www.naukri.com/job-listings-automation-testing-engineer-playwright-sonata-software-hyderabad-chennai-bengaluru-7-to-12-years-231225035819", "5 days ago", "Apply"),
    @("43.0", "Automation Test Engineer(C#+API)", "https://www.naukri.com/job-listings-automation-test-engineer-c-api-ltimindtree-hyderabad-chennai-bengaluru-5-to-10-years-020126008100", "5 days ago", "Apply"),
    @("44.0", "Test Automation Engineer", "https://www.naukri.com/job-listings-test-automation-engineer-accenture-solutions-pvt-ltd-chennai-3-to-8-years-130126921557", "5 days ago", "ApplyOnCompanySite"),
    @("45.0", "Storage Test Automation Engineer", "https://www.naukri.com/job-listings-storage-test-automation-engineer-aziro-pune-chennai-bengaluru-4-to-9-years-110126008534", "5 days ago", "Apply"),
    @("46.0", "Cypress Automation Test Engineer", "https://www.naukri.com/job-listings-cypress-automation-test-engineer-cognizant-chennai-4-to-10-years-020126013272", "5 days ago", "Apply"),
    @("47.0", "Test Automation Engineer | Contract", "https://www.naukri.com/job-listings-test-automation-engineer-contract-srs-infoway-chennai-8-to-10-years-090126939555", "5 days ago", "Apply"),
    @("48.0", "IVI Test Automation Engineer (Automotive Infotainment)", "https://www.naukri.com/job-listings-ivi-test-automation-engineer-automotive-infotainment-tekpillar-pune-chennai-bengaluru-5-to-10-years-100126000791", "5 days ago", "Apply"),
    @("49.0", "IVR Automation Test Engineer(Cyara & Kore.ai)", "https://www.naukri.com/job-listings-ivr-automation-test-engineer-cyara-kore-ai-summitworks-technologies-pune-chennai-5-to-10-years-030126013225", "5 days ago", "Apply"),
    @("50.0", "Automation Test Engineer", "https://www.naukri.com/job-listings-automation-test-engineer-indium-software-chennai-bengaluru-5-to-10-years-090126019771", "5 days ago", "Apply"),
    @("51.0", "Automation Test Engineer (UFT)", "https://www.naukri.com/job-listings-automation-test-engineer-uft-tata-consultancy-services-chennai-bengaluru-5-to-10-years-080126026105", "2 days ago", "Apply"),
    @("52.0", "Automation Testing Engineer- Playwright", "https://www.naukri.com/job-listings-automation-testing-engineer-playwright-sonata-software-hyderabad-chennai-bengaluru-7-to-12-years-231225035819", "2 days ago", "Apply"),
    @("53.0", "Automation Test Engineer(C#+API)", "https://www.naukri.com/job-listings-automation-test-engineer-c-api-ltimindtree-hyderabad-chennai-bengaluru-5-to-10-years-020126008100", "2 days ago", "Apply"),
    @("54.0", "Test Automation Engineer", "https://www.naukri.com/job-listings-test-automation-engineer-accenture-solutions-pvt-ltd-chennai-3-to-8-years-130126921557", "2 days ago", "ApplyOnCompanySite"),
    @("55.0", "Storage Test Automation Engineer", "https://www.naukri.com/job-listings-storage-test-automation-engineer-aziro-pune-chennai-bengaluru-4-to-9-years-110126008534", "2 days ago", "Apply"),
    @("56.0", "Cypress Automation Test Engineer", "https://www.naukri.com/job-listings-cypress-automation-test-engineer-cognizant-chennai-4-to-10-years-020126013272", "2 days ago", "Apply"),
    @("57.0", "Test Automation Engineer | Contract", "https://www.naukri.com/job-listings-test-automation-engineer-contract-srs-infoway-chennai-8-to-10-years-090126939555", "2 days ago", "Apply"),
    @("58.0", "IVI Test Automation Engineer (Automotive Infotainment)", "https://www.naukri.com/job-listings-ivi-test-automation-engineer-automotive-infotainment-tekpillar-pune-chennai-bengaluru-5-to-10-years-100126000791", "2 days ago", "Apply"),
    @("59.0", "IVR Automation Test Engineer(Cyara & Kore.ai)", "https://www.naukri.com/job-listings-ivr-automation-test-engineer-cyara-kore-ai-summitworks-technologies-pune-chennai-5-to-10-years-030126013225", "2 days ago", "Apply"),
    @("60.0", "Automation Test Engineer", "https://www.naukri.com/job-listings-automation-test-engineer-indium-software-chennai-bengaluru-5-to-10-years-090126019771", "2 days ago", "Apply"),
    @("61.0", "Automation Test Engineer (UFT)", "https://www.naukri.com/job-listings-automation-test-engineer-uft-tata-consultancy-services-chennai-bengaluru-5-to-10-years-080126026105", "1 day ago", "Apply"),
    @("62.0", "Automation Testing Engineer- Playwright", "https://www.naukri.com/job-listings-automation-testing-engineer-playwright-sonata-software-hyderabad-chennai-bengaluru-7-to-12-years-231225035819", "1 day ago", "Apply"),
    @("63.0", "Automation Test Engineer(C#+API)", "https://www.naukri.com/job-listings-automation-test-engineer-c-api-ltimindtree-hyderabad-chennai-bengaluru-5-to-10-years-020126008100", "1 day ago", "Apply"),
    @("64.0", "Test Automation Engineer", "https://www.naukri.com/job-listings-test-automation-engineer-accenture-solutions-pvt-ltd-chennai-3-to-8-years-130126921557", "1 day ago", "ApplyOnCompanySite"),
    @("65.0", "Storage Test Automation Engineer", "https://www.naukri.com/job-listings-storage-test-automation-engineer-aziro-pune-chennai-bengaluru-4-to-9-years-110126008534", "1 day ago", "Apply"),
    @("66.0", "Cypress Automation Test Engineer", "https://www.naukri.com/job-listings-cypress-automation-test-engineer-cognizant-chennai-4-to-10-years-020126013272", "1 day ago", "Apply"),
    @("67.0", "Test Automation Engineer | Contract", "https://www.naukri.com/job-listings-test-automation-engineer-contract-srs-infoway-chennai-8-to-10-years-090126939555", "1 day ago", "Apply"),
    @("68.0", "IVI Test Automation Engineer (Automotive Infotainment)", "https://www.naukri.com/job-listings-ivi-test-automation-engineer-automotive-infotainment-tekpillar-pune-chennai-bengaluru-5-to-10-years-100126000791", "1 day ago", "Apply"),
    @("69.0", "IVR Automation Test Engineer(Cyara & Kore.ai)", "https://www.naukri.com/job-listings-ivr-automation-test-engineer-cyara-kore-ai-summitworks-technologies-pune-chennai-5-to-10-years-030126013225", "1 day ago", "Apply"),
    @("70.0", "Automation Test Engineer", "https://www.naukri.com/job-listings-automation-test-engineer-indium-software-chennai-bengaluru-5-to-10-years-090126019771", "1 day ago", "Apply"),
    @("71.0", "Automation Test Engineer (UFT)", "https://www.naukri.com/job-listings-automation-test-engineer-uft-tata-consultancy-services-chennai-bengaluru-5-to-10-years-080126026105", "3 days ago", "Apply")
)

$startRow = 1412
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $r = $data[$i]
    $ws.Cells.Item($row, 1).Value = [double]$r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
}
